# Update HR and CI
# - Rename three header labels
# - Bold the header row
# - Resize a few columns to fit their (now longer / bold) header text
# - Restore the last-saved selection to F18

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header label changes -------------------------------------------------
# Set them in this order (I1, F1, C1) so that new shared-string entries are
# created in the same order Excel itself produced them.
$ws.Range("I1").Value = "mean_expectedHR"
$ws.Range("F1").Value = "datafrom"
$ws.Range("C1").Value = "interaction"

# --- Header row formatting --------------------------------------------------
# Bold the whole header row; existing alignment (center for most columns,
# left for column F) is preserved automatically.
$ws.Range("A1:I1").Font.Bold = $true

# --- Column width adjustments (to fit the new, now-bold, header text) ------
$ws.Columns.Item(1).ColumnWidth = 7.944010416666667
$ws.Columns.Item(2).ColumnWidth = 8.385416666666666
$ws.Columns.Item(3).ColumnWidth = 9.053385416666666
$ws.Columns.Item(4).ColumnWidth = 7.944010416666667
$ws.Columns.Item(5).ColumnWidth = 8.385416666666666
$ws.Columns.Item(6).ColumnWidth = 28.166666666666668
$ws.Columns.Item(7).ColumnWidth = 13.385416666666666
$ws.Columns.Item(8).ColumnWidth = 11.276041666666666
$ws.Columns.Item(9).ColumnWidth = 15.608072916666666

# --- Restore the active selection ------------------------------------------
$ws.Range("F18").Select()
